$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Checklist" to "Session"
$ws.Name = "Session"

# Update E2 value from "Selection" to "Scan"
$ws.Range("E2").Value = "Scan"

# Delete row 3 entirely (shifts nothing below it up, but removes its data/dimension)
$ws.Rows("3:3").Delete()
